$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 39
$ws.Range("I2").Value = 126
$ws.Range("J2").Value = 533
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 141
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 81
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 59
$ws.Range("T2").Value = 92
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 785
$ws.Range("X2").Value = 817
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 3
